# TablasTPTemporales.xlsx - edit per commit:
# "1_i hecha, maté tabla Licencia, taché pdf, modifiqué tabla Inasistencia,
#  corregí excel, CREATE e inserts de Docente en 1_a"
#
# Net effect on the worksheet (Hoja1):
#  - Drop the "Licencia" sub-table that used to sit at I14:M15.
#  - Move the "Inasistencia" sub-table from O14:R15 to J14:M15, renaming
#    its "Justifica" column header to "Justificación".
#  - Column M gets very slightly wider (new header text is longer).
#  - Selection/scroll position resets onto the moved table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the Inasistencia table's current contents before we move it ---
$inasistenciaTitle = $ws.Range("O14").Value()
$idInasistencia    = $ws.Range("O15").Value()
$fecha             = $ws.Range("P15").Value()
$cuil              = $ws.Range("Q15").Value()

# Unmerge every merged block we are about to touch so individual cells
# become writable (Excel refuses writes to non-anchor merged cells).
$ws.Range("I14:M14").UnMerge()
$ws.Range("O14:R14").UnMerge()

# --- wipe the whole old area (Licencia I14:M15) first ---
$ws.Range("I14:M15").Clear()

# --- copy formatting for the new J14:M14 / J15:M15 cells from the cells
#     that already carry the desired look (the old Inasistencia block),
#     AFTER the clear above so the freshly-pasted formats survive ---
$ws.Range("O14").Copy()
$ws.Range("J14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("O15").Copy()
$ws.Range("J15").PasteSpecial(-4122)

$ws.Range("P15").Copy()
$ws.Range("K15").PasteSpecial(-4122)

$ws.Range("Q15").Copy()
$ws.Range("L15").PasteSpecial(-4122)

$ws.Range("R15").Copy()
$ws.Range("M15").PasteSpecial(-4122)

# now it is safe to clear the old Inasistencia block (O14:R15) - it no
# longer overlaps anything we still need formatting from
$ws.Range("O14:R15").Clear()

# --- write the Inasistencia table into its new home, J14:M15 ---
$ws.Range("J14").Value = $inasistenciaTitle
$ws.Range("J15").Value = $idInasistencia
$ws.Range("K15").Value = $fecha
$ws.Range("L15").Value = $cuil
$ws.Range("M15").Value = "Justificación"

# re-merge the title cell over its new span
$ws.Range("J14:M14").Merge()

# column M needs to be a touch wider to fit "Justificación"
$ws.Columns.Item(13).ColumnWidth = 11

# reset scroll/selection onto the moved table
$ws.Range("O19").Select()
